# Casos de prueba actualizados
# Se definieron los valores a ingresar en los campos

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Precondiciones")
$ws3 = $wb.Worksheets.Item("Pasos")

# --- Precondiciones: set concrete values for the "CPA_" (Caso de Prueba A) scenario ---
$ws2.Range('B4').Value = '"CPA_Playa1" es el nombre de la playa <Playa1>'
$ws2.Range('B5').Value = ' "999999" es el telefono de la playa <Playa1>'
$ws2.Range('B6').Value = ' "CPA_mail@CPA_mail" es el mail de la playa <Playa1>'
$ws2.Range('B8').Value = ' <CPA_TipoVehiculo1> es un tipo de vehiculo aceptado por <Playa1> y su capacidad es <capacidad1>'
$ws2.Range('B7').Value = '<CPA_TipoPlaya1> es el tipo de playa de la playa <Playa1>'
$ws2.Range('B9').Value = '<Domicilio1>, conla provincia <Córdoba> existe en la base de datos,  el departamento <Capital> existe en la base de datos y pertenece a  <Córdoba>, la ciudad <Córdoba> existe en la base de datos y pertenece al <Capital>, <Calle1> con nombre "Colon" es una calle valida y "9" es un valor valido para numero; es el domicilio de la <Playa1>'
$ws2.Range('B10').Value = '<Domicilio2>, conla provincia <Córdoba> existe en la base de datos,  el departamento <Capital> existe en la base de datos y pertenece a  <Córdoba>, la ciudad <Córdoba> existe en la base de datos y pertenece al <Capital>, <Calle1> con nombre "Dean Funes" es una calle valida y "9" es un valor valido para numero; es domicilio valido para una playa de estacionamiento'
$ws2.Range('B11').Value = '<Horario1> con <CPA_DiasDeAtencion1> como dias, "00:00" como horario desde y "23:59" hasta es un horario de <Playa1>'
$ws2.Range('B12').Value = ' <Precio1> con <CPA_TipoVehiculo1> como tipo de vehiculo, <CPA_TipoHorario1> como tipo de horario, <CPA_DiasDeAtencion1> como dias y "9" como precio es un precio de <Playa1>'
$ws2.Range('B13').Value = '"CPA_Playa2" es un nombre válido para una playa de estcionamiento'

# --- Pasos: set concrete values entered in each field during the test ---
$ws3.Range('B3').Value = 'Ingreso "CPA_Playa1" en el campo nombre de playa'
$ws3.Range('B7').Value = 'Selecciono "Córdoba" como provincia'
$ws3.Range('B8').Value = 'Selecciono "Capital"  como departamento'
$ws3.Range('B9').Value = 'Selecciono "Córdoba" como ciudad'
$ws3.Range('B10').Value = 'Ingreso Dean Funes como calle'
$ws3.Range('B11').Value = 'Ingreso 9 como numero'
$ws3.Range('C13').Value = 'Se muestra el mensaje "Está seguro que desea guardar los cambios en la playa de estacionamiento CPA_Playa1?'

# --- View state: Precondiciones scrolled/selected near the bottom, Pasos becomes the active tab ---
$ws2.Range('B13').Select()

$ws3.Activate()
$ws3.Range('C14').Select()
